$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 66 (pushes the old blank separator + summary rows down
# by one, and Excel auto-extends the SUM(F2:F..) range to match).
$ws.Rows.Item(66).Insert()

# Fill in the new working-hours entry for 2014-03-13.
$ws.Cells.Item(66, 1).Value = 2014
$ws.Cells.Item(66, 2).Value = 3
$ws.Cells.Item(66, 3).Value = 13
$ws.Cells.Item(66, 4).Value = 0.78125
$ws.Cells.Item(66, 5).Value = 0.79861111111111116
$ws.Cells.Item(66, 6).Formula = "=(E66-D66)*24*60"
$ws.Cells.Item(66, 7).Formula = "=F66/60"

# Update the active selection to match the new last-used cell.
[void]$ws.Range("H69").Select()
